$p = $ppt.ActivePresentation

$oldDate = "3/18/2025"
$newDate = "3/19/2025"

function Update-DatePlaceholder {
    param($shapes)
    $n = $shapes.Count
    for ($i = 1; $i -le $n; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$m = $p.SlideMaster

# Update the date placeholder on the slide master itself.
Update-DatePlaceholder $m.Shapes

# Update the date placeholder on every slide layout (custom layout) under the master.
$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholder $lay.Shapes
}
